# C5-PowerPoint.pptx edit
#
# 1. Slide 6 has a table whose style is switched from the deck's custom
#    "Table_0" style to the built-in "Medium Style 2 - Accent 1" style.
# 2. The presentation's (slide master) theme colour palette is switched
#    from the "Integral" palette to the stock "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{546E536C-EE9F-422D-991B-65114A8685F8}")

# --- 2. Theme colours (Integral -> Office) ---------------------------
# Order returned by ThemeColorScheme.Item(n):
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2
#   7 accent3 8 accent4 9 accent5 10 accent6 11 hlink 12 folHlink
# RGB() style decimal = R + G*256 + B*65536 (i.e. hex BBGGRR)
$officeThemeRGB = @(
    0,            # dk1      000000
    16777215,     # lt1      FFFFFF
    6968388,      # dk2      44546A
    15132391,     # lt2      E7E6E6
    13998939,     # accent1  5B9BD5
    3243501,      # accent2  ED7D31
    10855845,     # accent3  A5A5A5
    49407,        # accent4  FFC000
    12874308,     # accent5  4472C4
    4697456,      # accent6  70AD47
    12673797,     # hlink    0563C1
    7491477       # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRGB[$i - 1]
}
